# Update from MV - monthly data refresh adds a new row (01-08-2021) and
# revises a handful of recent-month figures (rows 174-176).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 174 revisions ---
$ws.Range("B174").Value = 102414
$ws.Range("E174").Value = 36920
$ws.Range("F174").Value = 7796
$ws.Range("G174").Value = 4574
$ws.Range("R174").Value = 81645
$ws.Range("U174").Value = 29071
$ws.Range("V174").Value = 7027
$ws.Range("W174").Value = 3924

# --- Row 175 revisions ---
$ws.Range("B175").Value = 94993
$ws.Range("E175").Value = 29681
$ws.Range("F175").Value = 6882
$ws.Range("G175").Value = 5205
$ws.Range("H175").Value = 2347
$ws.Range("I175").Value = 2571
$ws.Range("R175").Value = 73519
$ws.Range("U175").Value = 23072
$ws.Range("V175").Value = 5767
$ws.Range("W175").Value = 4389
$ws.Range("X175").Value = 2073
$ws.Range("Y175").Value = 2174

# --- Row 176 revisions ---
$ws.Range("B176").Value = 88072
$ws.Range("D176").Value = 22216
$ws.Range("E176").Value = 31134
$ws.Range("F176").Value = 6704
$ws.Range("G176").Value = 4755
$ws.Range("H176").Value = 3458
$ws.Range("R176").Value = 69519
$ws.Range("T176").Value = 16524
$ws.Range("U176").Value = 24932
$ws.Range("V176").Value = 6038
$ws.Range("W176").Value = 4238
$ws.Range("X176").Value = 2557

# --- New row 177 (period 01-08-2021) ---
# Column A holds a text label that looks like a date ("01-08-2021"). A
# direct .Value assignment would be auto-converted to a date serial by
# Excel's input parser (and would stamp a number-format style onto the
# cell, which the source file doesn't have). Instead, build the literal
# text via a formula in a scratch cell, copy it, and paste-special just
# the values into A177 so it lands as plain shared-string text with the
# worksheet's default (unstyled) formatting - exactly like the rest of
# column A.
$ws.Range("Z1").Formula = "=""01-08-2021"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("A177").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1").ClearContents()

$ws.Range("B177").Value = 93932
$ws.Range("C177").Value = 16296
$ws.Range("D177").Value = 25976
$ws.Range("E177").Value = 31944
$ws.Range("F177").Value = 8741
$ws.Range("G177").Value = 6692
$ws.Range("H177").Value = 2084
$ws.Range("I177").Value = 2200
$ws.Range("J177").Value = 20022
$ws.Range("K177").Value = 3849
$ws.Range("L177").Value = 5834
$ws.Range("M177").Value = 7315
$ws.Range("N177").Value = 1368
$ws.Range("O177").Value = 893
$ws.Range("P177").Value = 438
$ws.Range("Q177").Value = 324
$ws.Range("R177").Value = 73909
$ws.Range("S177").Value = 12446
$ws.Range("T177").Value = 20141
$ws.Range("U177").Value = 24628
$ws.Range("V177").Value = 7373
$ws.Range("W177").Value = 5799
$ws.Range("X177").Value = 1646
$ws.Range("Y177").Value = 1875
